# contratos-8-2010.xlsx: fix formatting mangled while scraping floating point
# numbers. The scraper swapped Argentine-locale punctuation (dot thousands
# separator, comma decimal) for US-locale punctuation, but collided with plain
# text fields that happened to contain commas/periods too. Reproduce the exact
# same textual fix cell-by-cell; all target cells keep their original "Text"
# storage (they were never real numbers in the source workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Razon social" / "Nombre Fantasia" free-text cells where a stray comma got
# turned into a period by the same (mis-applied) locale fix.
$ws.Range("E43").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("F43").Value = "PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH"
$ws.Range("E48").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F48").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E50").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E81").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E83").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# "Importe" column. These cells are plain text holding es-AR formatted
# amounts (e.g. "1.070,00"), not real numbers. Setting .Value on a string
# that merely *looks* numeric would make Excel silently coerce it to a
# number and drop the formatting, so we briefly force Text number-formatting
# to preserve the literal, then restore the cell to its original ("Normal")
# style so no visible/structural formatting change is left behind.
$c = $ws.Range("H2")
$c.NumberFormat = "@"
$c.Value = "1070.00"
$c.Style = "Normal"
$c = $ws.Range("H3")
$c.NumberFormat = "@"
$c.Value = "410.00"
$c.Style = "Normal"
$c = $ws.Range("H4")
$c.NumberFormat = "@"
$c.Value = "3500.00"
$c.Style = "Normal"
$c = $ws.Range("H5")
$c.NumberFormat = "@"
$c.Value = "1775.00"
$c.Style = "Normal"
$c = $ws.Range("H6")
$c.NumberFormat = "@"
$c.Value = "770.00"
$c.Style = "Normal"
$c = $ws.Range("H7")
$c.NumberFormat = "@"
$c.Value = "209.50"
$c.Style = "Normal"
$c = $ws.Range("H8")
$c.NumberFormat = "@"
$c.Value = "59.00"
$c.Style = "Normal"
$c = $ws.Range("H9")
$c.NumberFormat = "@"
$c.Value = "2559.86"
$c.Style = "Normal"
$c = $ws.Range("H10")
$c.NumberFormat = "@"
$c.Value = "52010.00"
$c.Style = "Normal"
$c = $ws.Range("H11")
$c.NumberFormat = "@"
$c.Value = "85066.73"
$c.Style = "Normal"
$c = $ws.Range("H12")
$c.NumberFormat = "@"
$c.Value = "22467.83"
$c.Style = "Normal"
$c = $ws.Range("H13")
$c.NumberFormat = "@"
$c.Value = "6237.00"
$c.Style = "Normal"
$c = $ws.Range("H14")
$c.NumberFormat = "@"
$c.Value = "6996.50"
$c.Style = "Normal"
$c = $ws.Range("H15")
$c.NumberFormat = "@"
$c.Value = "1486.48"
$c.Style = "Normal"
$c = $ws.Range("H16")
$c.NumberFormat = "@"
$c.Value = "5131.50"
$c.Style = "Normal"
$c = $ws.Range("H17")
$c.NumberFormat = "@"
$c.Value = "6144.29"
$c.Style = "Normal"
$c = $ws.Range("H18")
$c.NumberFormat = "@"
$c.Value = "519.00"
$c.Style = "Normal"
$c = $ws.Range("H19")
$c.NumberFormat = "@"
$c.Value = "8217.96"
$c.Style = "Normal"
$c = $ws.Range("H20")
$c.NumberFormat = "@"
$c.Value = "475.00"
$c.Style = "Normal"
$c = $ws.Range("H21")
$c.NumberFormat = "@"
$c.Value = "1237.00"
$c.Style = "Normal"
$c = $ws.Range("H22")
$c.NumberFormat = "@"
$c.Value = "119.98"
$c.Style = "Normal"
$c = $ws.Range("H23")
$c.NumberFormat = "@"
$c.Value = "950.00"
$c.Style = "Normal"
$c = $ws.Range("H24")
$c.NumberFormat = "@"
$c.Value = "682.35"
$c.Style = "Normal"
$c = $ws.Range("H25")
$c.NumberFormat = "@"
$c.Value = "27.00"
$c.Style = "Normal"
$c = $ws.Range("H26")
$c.NumberFormat = "@"
$c.Value = "2220.00"
$c.Style = "Normal"
$c = $ws.Range("H27")
$c.NumberFormat = "@"
$c.Value = "411.40"
$c.Style = "Normal"
$c = $ws.Range("H28")
$c.NumberFormat = "@"
$c.Value = "1182.00"
$c.Style = "Normal"
$c = $ws.Range("H29")
$c.NumberFormat = "@"
$c.Value = "9314.14"
$c.Style = "Normal"
$c = $ws.Range("H30")
$c.NumberFormat = "@"
$c.Value = "555.00"
$c.Style = "Normal"
$c = $ws.Range("H31")
$c.NumberFormat = "@"
$c.Value = "927.33"
$c.Style = "Normal"
$c = $ws.Range("H32")
$c.NumberFormat = "@"
$c.Value = "8336.00"
$c.Style = "Normal"
$c = $ws.Range("H33")
$c.NumberFormat = "@"
$c.Value = "216.50"
$c.Style = "Normal"
$c = $ws.Range("H34")
$c.NumberFormat = "@"
$c.Value = "45.96"
$c.Style = "Normal"
$c = $ws.Range("H35")
$c.NumberFormat = "@"
$c.Value = "31.00"
$c.Style = "Normal"
$c = $ws.Range("H36")
$c.NumberFormat = "@"
$c.Value = "99.00"
$c.Style = "Normal"
$c = $ws.Range("H37")
$c.NumberFormat = "@"
$c.Value = "1139.00"
$c.Style = "Normal"
$c = $ws.Range("H38")
$c.NumberFormat = "@"
$c.Value = "18.40"
$c.Style = "Normal"
$c = $ws.Range("H39")
$c.NumberFormat = "@"
$c.Value = "3612.28"
$c.Style = "Normal"
$c = $ws.Range("H40")
$c.NumberFormat = "@"
$c.Value = "1290.00"
$c.Style = "Normal"
$c = $ws.Range("H41")
$c.NumberFormat = "@"
$c.Value = "1286.60"
$c.Style = "Normal"
$c = $ws.Range("H42")
$c.NumberFormat = "@"
$c.Value = "1275.00"
$c.Style = "Normal"
$c = $ws.Range("H43")
$c.NumberFormat = "@"
$c.Value = "90.00"
$c.Style = "Normal"
$c = $ws.Range("H44")
$c.NumberFormat = "@"
$c.Value = "12513.00"
$c.Style = "Normal"
$c = $ws.Range("H45")
$c.NumberFormat = "@"
$c.Value = "3906.00"
$c.Style = "Normal"
$c = $ws.Range("H46")
$c.NumberFormat = "@"
$c.Value = "1200.00"
$c.Style = "Normal"
$c = $ws.Range("H47")
$c.NumberFormat = "@"
$c.Value = "483.00"
$c.Style = "Normal"
$c = $ws.Range("H48")
$c.NumberFormat = "@"
$c.Value = "5208.97"
$c.Style = "Normal"
$c = $ws.Range("H49")
$c.NumberFormat = "@"
$c.Value = "1164.12"
$c.Style = "Normal"
$c = $ws.Range("H50")
$c.NumberFormat = "@"
$c.Value = "320.00"
$c.Style = "Normal"
$c = $ws.Range("H51")
$c.NumberFormat = "@"
$c.Value = "200340.00"
$c.Style = "Normal"
$c = $ws.Range("H52")
$c.NumberFormat = "@"
$c.Value = "600.00"
$c.Style = "Normal"
$c = $ws.Range("H53")
$c.NumberFormat = "@"
$c.Value = "6000.00"
$c.Style = "Normal"
$c = $ws.Range("H54")
$c.NumberFormat = "@"
$c.Value = "5928.00"
$c.Style = "Normal"
$c = $ws.Range("H55")
$c.NumberFormat = "@"
$c.Value = "4.80"
$c.Style = "Normal"
$c = $ws.Range("H56")
$c.NumberFormat = "@"
$c.Value = "908.50"
$c.Style = "Normal"
$c = $ws.Range("H57")
$c.NumberFormat = "@"
$c.Value = "5700.30"
$c.Style = "Normal"
$c = $ws.Range("H58")
$c.NumberFormat = "@"
$c.Value = "102.72"
$c.Style = "Normal"
$c = $ws.Range("H59")
$c.NumberFormat = "@"
$c.Value = "825.00"
$c.Style = "Normal"
$c = $ws.Range("H60")
$c.NumberFormat = "@"
$c.Value = "13300.00"
$c.Style = "Normal"
$c = $ws.Range("H61")
$c.NumberFormat = "@"
$c.Value = "539.00"
$c.Style = "Normal"
$c = $ws.Range("H62")
$c.NumberFormat = "@"
$c.Value = "450.00"
$c.Style = "Normal"
$c = $ws.Range("H63")
$c.NumberFormat = "@"
$c.Value = "1800.00"
$c.Style = "Normal"
$c = $ws.Range("H64")
$c.NumberFormat = "@"
$c.Value = "970.00"
$c.Style = "Normal"
$c = $ws.Range("H65")
$c.NumberFormat = "@"
$c.Value = "479.95"
$c.Style = "Normal"
$c = $ws.Range("H66")
$c.NumberFormat = "@"
$c.Value = "39.80"
$c.Style = "Normal"
$c = $ws.Range("H67")
$c.NumberFormat = "@"
$c.Value = "2872.88"
$c.Style = "Normal"
$c = $ws.Range("H68")
$c.NumberFormat = "@"
$c.Value = "250.00"
$c.Style = "Normal"
$c = $ws.Range("H69")
$c.NumberFormat = "@"
$c.Value = "500.00"
$c.Style = "Normal"
$c = $ws.Range("H70")
$c.NumberFormat = "@"
$c.Value = "3213.76"
$c.Style = "Normal"
$c = $ws.Range("H71")
$c.NumberFormat = "@"
$c.Value = "290.00"
$c.Style = "Normal"
$c = $ws.Range("H72")
$c.NumberFormat = "@"
$c.Value = "500.00"
$c.Style = "Normal"
$c = $ws.Range("H73")
$c.NumberFormat = "@"
$c.Value = "1210.00"
$c.Style = "Normal"
$c = $ws.Range("H74")
$c.NumberFormat = "@"
$c.Value = "1379.04"
$c.Style = "Normal"
$c = $ws.Range("H75")
$c.NumberFormat = "@"
$c.Value = "200.00"
$c.Style = "Normal"
$c = $ws.Range("H76")
$c.NumberFormat = "@"
$c.Value = "350.00"
$c.Style = "Normal"
$c = $ws.Range("H77")
$c.NumberFormat = "@"
$c.Value = "750.00"
$c.Style = "Normal"
$c = $ws.Range("H78")
$c.NumberFormat = "@"
$c.Value = "120.00"
$c.Style = "Normal"
$c = $ws.Range("H79")
$c.NumberFormat = "@"
$c.Value = "3928.24"
$c.Style = "Normal"
$c = $ws.Range("H80")
$c.NumberFormat = "@"
$c.Value = "319.00"
$c.Style = "Normal"
$c = $ws.Range("H81")
$c.NumberFormat = "@"
$c.Value = "45.00"
$c.Style = "Normal"
$c = $ws.Range("H82")
$c.NumberFormat = "@"
$c.Value = "665.00"
$c.Style = "Normal"
$c = $ws.Range("H83")
$c.NumberFormat = "@"
$c.Value = "40.00"
$c.Style = "Normal"
$c = $ws.Range("H84")
$c.NumberFormat = "@"
$c.Value = "580.00"
$c.Style = "Normal"
$c = $ws.Range("H85")
$c.NumberFormat = "@"
$c.Value = "1150.00"
$c.Style = "Normal"
$c = $ws.Range("H86")
$c.NumberFormat = "@"
$c.Value = "22000.00"
$c.Style = "Normal"
$c = $ws.Range("H87")
$c.NumberFormat = "@"
$c.Value = "360.40"
$c.Style = "Normal"
$c = $ws.Range("H88")
$c.NumberFormat = "@"
$c.Value = "101.70"
$c.Style = "Normal"
$c = $ws.Range("H89")
$c.NumberFormat = "@"
$c.Value = "117.00"
$c.Style = "Normal"
$c = $ws.Range("H90")
$c.NumberFormat = "@"
$c.Value = "36.00"
$c.Style = "Normal"
$c = $ws.Range("H91")
$c.NumberFormat = "@"
$c.Value = "520.96"
$c.Style = "Normal"
$c = $ws.Range("H92")
$c.NumberFormat = "@"
$c.Value = "49.90"
$c.Style = "Normal"
$c = $ws.Range("H93")
$c.NumberFormat = "@"
$c.Value = "847.00"
$c.Style = "Normal"
$c = $ws.Range("H94")
$c.NumberFormat = "@"
$c.Value = "2688.88"
$c.Style = "Normal"
$c = $ws.Range("H95")
$c.NumberFormat = "@"
$c.Value = "2341.35"
$c.Style = "Normal"
$c = $ws.Range("H96")
$c.NumberFormat = "@"
$c.Value = "968.00"
$c.Style = "Normal"
$c = $ws.Range("H97")
$c.NumberFormat = "@"
$c.Value = "9958.48"
$c.Style = "Normal"
$c = $ws.Range("H98")
$c.NumberFormat = "@"
$c.Value = "2957.28"
$c.Style = "Normal"
$c = $ws.Range("H99")
$c.NumberFormat = "@"
$c.Value = "936896.62"
$c.Style = "Normal"
$c = $ws.Range("H100")
$c.NumberFormat = "@"
$c.Value = "44000.00"
$c.Style = "Normal"
$c = $ws.Range("H101")
$c.NumberFormat = "@"
$c.Value = "5700.00"
$c.Style = "Normal"
